$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Better LINQ" column from the results grid (column D).
#    This shifts MapReduce/Hash/Character Sequence left (E->D, F->E, G->F)
#    and also shifts the detail list in column I left into column H.
$ws.Columns.Item(4).Delete()

# 2. Consolidate the "LINQ" and "Better LINQ" detail rows in column H.
#    The "Better LINQ" row in each group keeps its numeric value but is
#    relabelled "LINQ"; the old (now redundant) "LINQ" row is cleared.
$ws.Range("H6").Value2  = "Search for a 10000 times using LINQ took 2754213400ns,2754.2134ms,2.7542134s"
$ws.Range("H7").ClearContents()

$ws.Range("H11").Value2 = "Search for PHOENIX 10000 times using LINQ took 1351685100ns,1351.6851ms,1.3516851s"
$ws.Range("H12").ClearContents()

$ws.Range("H16").Value2 = "Search for phoenix 10000 times using LINQ took 1296197700ns,1296.1977ms,1.2961977s"
$ws.Range("H17").ClearContents()

$ws.Range("H21").Value2 = "Search for catherine 10000 times using LINQ took 1155302400ns,1155.3024ms,1.1553024s"
$ws.Range("H22").ClearContents()

$ws.Range("H26").Value2 = "LINQ indexed 10000 took 39337700ns,39.3377ms,0.0393377s"
$ws.Range("H27").ClearContents()

# 3. Match the final selection state left behind by the edit.
$ws.Range("H26").Select()
